# Auto-generated edit script applying value updates per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 290.81818
$ws.Range("I8").Value = 290.81818
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 872.45454
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -733.45454
$ws.Range("N8").ClearContents()
$ws.Range("H9").Value = 523.8570999999999
$ws.Range("I9").Value = 472.6
$ws.Range("K9").Value = 472.6
$ws.Range("M9").Value = -303.6
$ws.Range("H15").Value = 4769.7354
$ws.Range("I15").Value = 4769.7354
$ws.Range("K15").Value = 14309.2062
$ws.Range("M15").Value = -14140.2062
$ws.Range("H17").Value = 1256.0217
$ws.Range("J17").Value = 1239.4889
$ws.Range("L17").Value = 3718.4667
$ws.Range("N17").Value = -4054.4667
$ws.Range("H33").Value = 912.7917
$ws.Range("I33").Value = 1087.3334
$ws.Range("K33").Value = 1087.3334
$ws.Range("M33").Value = -858.3334
$ws.Range("H53").Value = 462.5
$ws.Range("I53").Value = 491.2
$ws.Range("J53").Value = 433.8
$ws.Range("K53").Value = 491.2
$ws.Range("L53").Value = 433.8
$ws.Range("M53").Value = 145.8
$ws.Range("N53").Value = -1707.8
$ws.Range("H62").Value = 7749
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7749
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
$ws.Range("H76").Value = 4165.364
$ws.Range("J76").Value = 4290.2
$ws.Range("L76").Value = 4290.2
$ws.Range("N76").Value = -4920.2
$ws.Range("H79").Value = 4165.364
$ws.Range("J79").Value = 4290.2
$ws.Range("L79").Value = 4290.2
$ws.Range("N79").Value = -6474.2
$ws.Range("H98").Value = 1260.3334
$ws.Range("J98").Value = 493
$ws.Range("L98").Value = 493
$ws.Range("N98").Value = -3489
$ws.Range("H111").Value = 1765.2
$ws.Range("I111").Value = 1765.2
$ws.Range("K111").Value = 5295.6
$ws.Range("M111").Value = -2228.6
$ws.Range("H112").Value = 2277.7173
$ws.Range("I112").Value = 3649.3333
$ws.Range("J112").Value = 2071.975
$ws.Range("K112").Value = 10947.9999
$ws.Range("L112").Value = 6215.924999999999
$ws.Range("M112").Value = -9839.999899999999
$ws.Range("N112").Value = -8431.924999999999
$ws.Range("H120").Value = 70000
$ws.Range("J120").Value = 70000
$ws.Range("L120").Value = 70000
$ws.Range("N120").Value = -79676
$ws.Range("H121").Value = 7533
$ws.Range("J121").Value = 7249.5
$ws.Range("L121").Value = 21748.5
$ws.Range("N121").Value = -25242.5
$ws.Range("H122").Value = 1260.3334
$ws.Range("J122").Value = 493
$ws.Range("L122").Value = 1479
$ws.Range("N122").Value = -6379
$ws.Range("H128").Value = 200000
$ws.Range("J128").Value = 200000
$ws.Range("L128").Value = 200000
$ws.Range("N128").Value = -209960
$ws.Range("H131").Value = 7756.357
$ws.Range("J131").Value = 8249.875
$ws.Range("L131").Value = 24749.625
$ws.Range("N131").Value = -34829.625
$ws.Range("H132").Value = 10945.536
$ws.Range("I132").Value = 6403.5
$ws.Range("J132").Value = 19121.2
$ws.Range("K132").Value = 19210.5
$ws.Range("L132").Value = 57363.60000000001
$ws.Range("M132").Value = -16680.5
$ws.Range("N132").Value = -62423.60000000001
$ws.Range("H136").Value = 164848.6
$ws.Range("J136").Value = 164848.6
$ws.Range("L136").Value = 164848.6
$ws.Range("N136").Value = -175048.6
$ws.Range("H137").Value = 3478.7544
$ws.Range("I137").Value = 4032.5715
$ws.Range("J137").Value = 1928.0667
$ws.Range("K137").Value = 12097.7145
$ws.Range("L137").Value = 5784.2001
$ws.Range("M137").Value = -9547.7145
$ws.Range("N137").Value = -10884.2001
$ws.Range("H138").Value = 16669105
$ws.Range("I138").Value = 55556932
$ws.Range("K138").Value = 166670796
$ws.Range("M138").Value = -166665656
$ws.Range("H141").Value = 10888.565
$ws.Range("I141").Value = 11269.863
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 33809.589
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = -28629.589
$ws.Range("N141").Value = -17860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 59949
$ws.Range("J7").Value = 59949
$ws.Range("L7").Value = 59949
$ws.Range("N7").Value = -60177
$ws.Range("H32").Value = 4664.9663
$ws.Range("I32").Value = 4664.9663
$ws.Range("K32").Value = 4664.9663
$ws.Range("M32").Value = -4377.9663
$ws.Range("H45").Value = 8771.727999999999
$ws.Range("I45").Value = 11954.615
$ws.Range("J45").Value = 4174.222
$ws.Range("K45").Value = 11954.615
$ws.Range("L45").Value = 4174.222
$ws.Range("M45").Value = -11577.615
$ws.Range("N45").Value = -4928.222
$ws.Range("H61").Value = 5009.2856
$ws.Range("I61").Value = 4980.147
$ws.Range("K61").Value = 4980.147
$ws.Range("M61").Value = -4768.147
$ws.Range("H120").Value = 112920.336
$ws.Range("J120").Value = 112920.336
$ws.Range("L120").Value = 112920.336
$ws.Range("N120").Value = -122596.336
$ws.Range("H121").Value = 172965
$ws.Range("J121").Value = 172965
$ws.Range("L121").Value = 172965
$ws.Range("N121").Value = -176459
$ws.Range("H122").Value = 17991.25
$ws.Range("I122").Value = 1604.2941
$ws.Range("J122").Value = 110850.664
$ws.Range("K122").Value = 4812.8823
$ws.Range("L122").Value = 332551.992
$ws.Range("M122").Value = -2362.8823
$ws.Range("N122").Value = -337451.992
$ws.Range("H123").Value = 49995
$ws.Range("J123").Value = 49995
$ws.Range("L123").Value = 49995
$ws.Range("N123").Value = -59795
$ws.Range("H132").Value = 2432.087
$ws.Range("I132").Value = 2147.05
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 6441.150000000001
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -3911.150000000001
$ws.Range("N132").Value = -18057.0005
$ws.Range("H134").Value = 70000
$ws.Range("J134").Value = 70000
$ws.Range("L134").Value = 70000
$ws.Range("N134").Value = -80140
$ws.Range("H136").Value = 5009.2856
$ws.Range("I136").Value = 4980.147
$ws.Range("K136").Value = 14940.441
$ws.Range("M136").Value = -12390.441
$ws.Range("H140").Value = 132621.25
$ws.Range("I140").Value = 153494
$ws.Range("J140").Value = 125663.664
$ws.Range("K140").Value = 153494
$ws.Range("L140").Value = 125663.664
$ws.Range("M140").Value = -148314
$ws.Range("N140").Value = -136023.664

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6468.467
$ws.Range("I20").Value = 5419.7
$ws.Range("J20").Value = 8566
$ws.Range("K20").Value = 5419.7
$ws.Range("L20").Value = 8566
$ws.Range("M20").Value = -5172.7
$ws.Range("N20").Value = -9060
$ws.Range("H38").Value = 16024
$ws.Range("J38").Value = 16024
$ws.Range("L38").Value = 16024
$ws.Range("N38").Value = -16950
$ws.Range("H102").Value = 38171.234
$ws.Range("I102").Value = 47363.875
$ws.Range("K102").Value = 47363.875
$ws.Range("M102").Value = -44118.875
$ws.Range("H107").Value = 4445.3
$ws.Range("I107").Value = 3661.4119
$ws.Range("J107").Value = 8887.333000000001
$ws.Range("K107").Value = 3661.4119
$ws.Range("L107").Value = 8887.333000000001
$ws.Range("M107").Value = -1741.4119
$ws.Range("N107").Value = -12727.333
$ws.Range("H134").Value = 5426
$ws.Range("I134").Value = 5426
$ws.Range("K134").Value = 16278
$ws.Range("M134").Value = -13743

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2173
$ws.Range("I31").Value = 2094.9412
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 2094.9412
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = -1799.9412
$ws.Range("N31").Value = -4090
$ws.Range("H34").Value = 2173
$ws.Range("I34").Value = 2094.9412
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 2094.9412
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = -1892.9412
$ws.Range("N34").Value = -3904
$ws.Range("H58").Value = 3476.3333
$ws.Range("I58").Value = 3445
$ws.Range("J58").Value = 3485.2856
$ws.Range("K58").Value = 3445
$ws.Range("L58").Value = 3485.2856
$ws.Range("M58").Value = -3242
$ws.Range("N58").Value = -3891.2856
$ws.Range("H94").Value = 8887
$ws.Range("I94").Value = 8887
$ws.Range("K94").Value = 8887
$ws.Range("M94").Value = -8436
$ws.Range("H100").Value = 66325
$ws.Range("J100").Value = 66325
$ws.Range("L100").Value = 66325
$ws.Range("N100").Value = -68489
$ws.Range("H119").Value = 120000
$ws.Range("I119").Value = 100000
$ws.Range("J119").Value = 140000
$ws.Range("K119").Value = 100000
$ws.Range("L119").Value = 140000
$ws.Range("M119").Value = -95162
$ws.Range("N119").Value = -149676
$ws.Range("H122").Value = 3193.889
$ws.Range("I122").Value = 2669.3635
$ws.Range("K122").Value = 8008.0905
$ws.Range("M122").Value = -5558.0905
$ws.Range("H132").Value = 21893.46
$ws.Range("I132").Value = 8946.286
$ws.Range("J132").Value = 36998.5
$ws.Range("K132").Value = 26838.858
$ws.Range("L132").Value = 110995.5
$ws.Range("M132").Value = -24308.858
$ws.Range("N132").Value = -116055.5
$ws.Range("H134").Value = 3439.375
$ws.Range("I134").Value = 3168.8572
$ws.Range("K134").Value = 9506.571599999999
$ws.Range("M134").Value = -6971.571599999999
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140
$ws.Range("H136").Value = 3476.3333
$ws.Range("I136").Value = 3445
$ws.Range("J136").Value = 3485.2856
$ws.Range("K136").Value = 10335
$ws.Range("L136").Value = 10455.8568
$ws.Range("M136").Value = -7785
$ws.Range("N136").Value = -15555.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 361.77777
$ws.Range("I5").Value = 294.82352
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 884.47056
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -772.47056
$ws.Range("N5").Value = -4724
$ws.Range("H36").Value = 596
$ws.Range("I36").Value = 596
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1788
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1619
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 3951
$ws.Range("J46").Value = 2276.1667
$ws.Range("L46").Value = 6828.500100000001
$ws.Range("N46").Value = -7010.500100000001
$ws.Range("H62").Value = 15999.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 15999.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 47998.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -49370.5
$ws.Range("H65").Value = 15999.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 15999.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 143995.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -150859.5
$ws.Range("H74").Value = 4998.5
$ws.Range("I74").Value = 4998.5
$ws.Range("K74").Value = 14995.5
$ws.Range("M74").Value = -13934.5
$ws.Range("H77").Value = 4998.5
$ws.Range("I77").Value = 4998.5
$ws.Range("K77").Value = 44986.5
$ws.Range("M77").Value = -39682.5
$ws.Range("H87").Value = 4749.5
$ws.Range("I87").Value = 4749
$ws.Range("K87").Value = 14247
$ws.Range("M87").Value = -12999
$ws.Range("H90").Value = 4749.5
$ws.Range("I90").Value = 4749
$ws.Range("K90").Value = 42741
$ws.Range("M90").Value = -36501
$ws.Range("H98").Value = 596.8889
$ws.Range("I98").Value = 620.4
$ws.Range("J98").Value = 567.5
$ws.Range("K98").Value = 1861.2
$ws.Range("L98").Value = 1702.5
$ws.Range("M98").Value = -363.1999999999998
$ws.Range("N98").Value = -4698.5
$ws.Range("H107").Value = 562.6
$ws.Range("I107").Value = 504.57144
$ws.Range("K107").Value = 1513.71432
$ws.Range("M107").Value = 406.28568
$ws.Range("H113").Value = 2309.0527
$ws.Range("I113").Value = 2301.6667
$ws.Range("J113").Value = 2336.75
$ws.Range("K113").Value = 6905.000100000001
$ws.Range("L113").Value = 7010.25
$ws.Range("M113").Value = -4735.000100000001
$ws.Range("N113").Value = -11350.25
$ws.Range("H116").Value = 2544.2
$ws.Range("I116").Value = 2582
$ws.Range("J116").Value = 2487.5
$ws.Range("K116").Value = 7746
$ws.Range("L116").Value = 7462.5
$ws.Range("M116").Value = -4304
$ws.Range("N116").Value = -14346.5
$ws.Range("H119").Value = 3443.2856
$ws.Range("I119").Value = 2361.8333
$ws.Range("K119").Value = 7085.499899999999
$ws.Range("M119").Value = -2247.499899999999
$ws.Range("H122").Value = 897.8889
$ws.Range("I122").Value = 508.69232
$ws.Range("J122").Value = 1259.2858
$ws.Range("K122").Value = 4578.23088
$ws.Range("L122").Value = 11333.5722
$ws.Range("M122").Value = -2128.23088
$ws.Range("N122").Value = -16233.5722
$ws.Range("H124").Value = 8253.333000000001
$ws.Range("J124").Value = 12500
$ws.Range("L124").Value = 37500
$ws.Range("N124").Value = -47320
$ws.Range("H125").Value = 9932.25
$ws.Range("H135").Value = 361.77777
$ws.Range("I135").Value = 294.82352
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 2653.41168
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -118.4116799999997
$ws.Range("N135").Value = -18570

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70707150
$ws.Range("I2").Value = 111111170
$ws.Range("J2").Value = 113.5
$ws.Range("K2").Value = 111111170
$ws.Range("L2").Value = 113.5
$ws.Range("M2").Value = -111111057
$ws.Range("N2").Value = -339.5
$ws.Range("H22").Value = 4999.3335
$ws.Range("I22").Value = 2999
$ws.Range("J22").Value = 5999.5
$ws.Range("K22").Value = 2999
$ws.Range("L22").Value = 5999.5
$ws.Range("M22").Value = -2470
$ws.Range("N22").Value = -7057.5
$ws.Range("H25").Value = 2858.25
$ws.Range("I25").Value = 3100
$ws.Range("J25").Value = 2616.5
$ws.Range("K25").Value = 3100
$ws.Range("L25").Value = 2616.5
$ws.Range("M25").Value = -2571
$ws.Range("N25").Value = -3674.5
$ws.Range("H70").Value = 7349
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 7349
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7349
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -7889
$ws.Range("H73").Value = 7349
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 7349
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7349
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9221
$ws.Range("H97").Value = 545.75
$ws.Range("I97").Value = 545.75
$ws.Range("K97").Value = 545.75
$ws.Range("M97").Value = -49.75
$ws.Range("H102").Value = 6982.722
$ws.Range("I102").Value = 6911.7334
$ws.Range("K102").Value = 6911.7334
$ws.Range("M102").Value = -5289.7334
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 2763.0356
$ws.Range("I122").Value = 2268.8235
$ws.Range("K122").Value = 6806.470499999999
$ws.Range("M122").Value = -4356.470499999999
$ws.Range("H132").Value = 3581.8823
$ws.Range("I132").Value = 1724.1111
$ws.Range("K132").Value = 5172.3333
$ws.Range("M132").Value = -2642.3333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1356.9231
$ws.Range("I7").Value = 1386.8334
$ws.Range("J7").Value = 998
$ws.Range("K7").Value = 1386.8334
$ws.Range("L7").Value = 998
$ws.Range("M7").Value = -1274.8334
$ws.Range("N7").Value = -1222
$ws.Range("H22").Value = 667.1111
$ws.Range("I22").Value = 626.125
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 626.125
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -331.125
$ws.Range("N22").Value = -1585
$ws.Range("H27").Value = 667.1111
$ws.Range("I27").Value = 626.125
$ws.Range("J27").Value = 995
$ws.Range("K27").Value = 626.125
$ws.Range("L27").Value = 995
$ws.Range("M27").Value = -519.125
$ws.Range("N27").Value = -1209
$ws.Range("H40").Value = 4584.2144
$ws.Range("I40").Value = 4480.8184
$ws.Range("J40").Value = 4963.3335
$ws.Range("K40").Value = 4480.8184
$ws.Range("L40").Value = 4963.3335
$ws.Range("M40").Value = -4344.8184
$ws.Range("N40").Value = -5235.3335
$ws.Range("H46").Value = 2478.7
$ws.Range("I46").Value = 999.8
$ws.Range("J46").Value = 3218.15
$ws.Range("K46").Value = 999.8
$ws.Range("L46").Value = 3218.15
$ws.Range("M46").Value = -811.8
$ws.Range("N46").Value = -3594.15
$ws.Range("H55").Value = 198.30435
$ws.Range("I55").Value = 138.1
$ws.Range("J55").Value = 599.6667
$ws.Range("K55").Value = 138.1
$ws.Range("L55").Value = 599.6667
$ws.Range("M55").Value = 34.90000000000001
$ws.Range("N55").Value = -945.6667
$ws.Range("H82").Value = 2233.61
$ws.Range("I82").Value = 2269.5833
$ws.Range("J82").Value = 1370.25
$ws.Range("K82").Value = 2269.5833
$ws.Range("L82").Value = 1370.25
$ws.Range("M82").Value = -1908.5833
$ws.Range("N82").Value = -2092.25
$ws.Range("H85").Value = 2233.61
$ws.Range("I85").Value = 2269.5833
$ws.Range("J85").Value = 1370.25
$ws.Range("K85").Value = 2269.5833
$ws.Range("L85").Value = 1370.25
$ws.Range("M85").Value = -1021.5833
$ws.Range("N85").Value = -3866.25
$ws.Range("H122").Value = 5953.2383
$ws.Range("I122").Value = 6529.857
$ws.Range("K122").Value = 19589.571
$ws.Range("M122").Value = -17139.571
$ws.Range("H126").Value = 1356.9231
$ws.Range("I126").Value = 1386.8334
$ws.Range("J126").Value = 998
$ws.Range("K126").Value = 4160.5002
$ws.Range("L126").Value = 2994
$ws.Range("M126").Value = -1690.5002
$ws.Range("N126").Value = -7934
$ws.Range("H130").Value = 148199.25
$ws.Range("J130").Value = 148199.25
$ws.Range("L130").Value = 148199.25
$ws.Range("N130").Value = -158239.25
$ws.Range("H132").Value = 3916.3157
$ws.Range("I132").Value = 3856.1667
$ws.Range("K132").Value = 11568.5001
$ws.Range("M132").Value = -9038.500100000001
$ws.Range("H134").Value = 113778
$ws.Range("J134").Value = 113778
$ws.Range("L134").Value = 113778
$ws.Range("N134").Value = -123918
$ws.Range("H135").Value = 71214.336
$ws.Range("J135").Value = 71214.336
$ws.Range("L135").Value = 71214.336
$ws.Range("N135").Value = -81354.336
$ws.Range("H136").Value = 3298.8604
$ws.Range("I136").Value = 2663.3333
$ws.Range("J136").Value = 3467.0881
$ws.Range("K136").Value = 7989.999899999999
$ws.Range("L136").Value = 10401.2643
$ws.Range("M136").Value = -5439.999899999999
$ws.Range("N136").Value = -15501.2643
$ws.Range("H141").Value = 81464
$ws.Range("J141").Value = 81464
$ws.Range("L141").Value = 81464
$ws.Range("N141").Value = -91824

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 31721.777
$ws.Range("J40").Value = 26500
$ws.Range("L40").Value = 26500
$ws.Range("N40").Value = -26798
$ws.Range("H45").Value = 13067.3
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H54").Value = 35099.285
$ws.Range("I54").Value = 40000
$ws.Range("J54").Value = 34282.5
$ws.Range("K54").Value = 40000
$ws.Range("L54").Value = 34282.5
$ws.Range("M54").Value = -39480
$ws.Range("N54").Value = -35322.5
$ws.Range("H94").Value = 44666
$ws.Range("J94").Value = 44666
$ws.Range("L94").Value = 44666
$ws.Range("N94").Value = -46468
$ws.Range("H113").Value = 846.02325
$ws.Range("I113").Value = 729.40625
$ws.Range("J113").Value = 1185.2727
$ws.Range("K113").Value = 2188.21875
$ws.Range("L113").Value = 3555.8181
$ws.Range("M113").Value = -18.21875
$ws.Range("N113").Value = -7895.8181
$ws.Range("H122").Value = 3714.6416
$ws.Range("I122").Value = 1875.1555
$ws.Range("K122").Value = 5625.4665
$ws.Range("M122").Value = -3175.4665
$ws.Range("H127").Value = 77100
$ws.Range("J127").Value = 77100
$ws.Range("L127").Value = 77100
$ws.Range("N127").Value = -87020
$ws.Range("H131").Value = 86857.5
$ws.Range("J131").Value = 86857.5
$ws.Range("L131").Value = 86857.5
$ws.Range("N131").Value = -96937.5
$ws.Range("H132").Value = 3756.8
$ws.Range("I132").Value = 3017.92
$ws.Range("J132").Value = 4988.2666
$ws.Range("K132").Value = 9053.76
$ws.Range("L132").Value = 14964.7998
$ws.Range("M132").Value = -6523.76
$ws.Range("N132").Value = -20024.7998
$ws.Range("H136").Value = 6291
$ws.Range("I136").Value = 6338.6665
$ws.Range("K136").Value = 19015.9995
$ws.Range("M136").Value = -16465.9995
$ws.Range("H137").Value = 114748.2
$ws.Range("J137").Value = 114748.2
$ws.Range("L137").Value = 114748.2
$ws.Range("N137").Value = -124948.2
$ws.Range("H138").Value = 55000
$ws.Range("J138").Value = 55000
$ws.Range("L138").Value = 55000
$ws.Range("N138").Value = -65280
$ws.Range("H139").Value = 69913.57000000001
$ws.Range("J139").Value = 69913.57000000001
$ws.Range("L139").Value = 69913.57000000001
$ws.Range("N139").Value = -80193.57000000001
$ws.Range("H140").Value = 66043.89
$ws.Range("J140").Value = 66043.89
$ws.Range("L140").Value = 66043.89
$ws.Range("N140").Value = -76403.89
